$d = $word.ActiveDocument

# Helper: scan all paragraphs and return the LAST one whose style name and
# text match (the body headings come after the generated Table of Contents
# entries, which use "TOC n" paragraph styles rather than "Heading n", so a
# plain last-match lookup unambiguously lands on the real body heading).
function Find-ParaByStyleLike($styleName, $pattern) {
    $result = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Style.NameLocal -eq $styleName -and $p.Range.Text -like $pattern) {
            $result = $p
        }
    }
    return $result
}

# ---------------------------------------------------------------------------
# 1) Remove the whole "Chapitre B - Gestion des utilisateurs" playground
#    sub-section in the "Spécificités" part: the Titre3 heading plus its two
#    bullet paragraphs plus the blank spacer paragraph that follows it, i.e.
#    everything from that heading up to (not including) "Chapitre F".
# ---------------------------------------------------------------------------
$chapB = Find-ParaByStyleLike "Heading 3" "*Gestion des utilisateurs*"
$chapF = Find-ParaByStyleLike "Heading 3" "*Chapitre F*"
if ($chapB -eq $null -or $chapF -eq $null) {
    Write-Output "ERROR: could not locate Chapitre B / Chapitre F headings"
} else {
    $d.Range($chapB.Range.Start, $chapF.Range.Start).Delete()
}

# ---------------------------------------------------------------------------
# 2) In the "Pratique" chapter overview (the block of Titre2 headings), drop
#    the "B", "C", "D" and "E" placeholder headings - only "A", "F" and "G"
#    remain.
# ---------------------------------------------------------------------------
$pratiqueHeadB = Find-ParaByStyleLike "Heading 2" "*Gestion des utilisateurs*"
$pratiqueHeadF = Find-ParaByStyleLike "Heading 2" "*Création des index*"
if ($pratiqueHeadB -eq $null -or $pratiqueHeadF -eq $null) {
    Write-Output "ERROR: could not locate B / F overview headings"
} else {
    $d.Range($pratiqueHeadB.Range.Start, $pratiqueHeadF.Range.Start).Delete()
}

# ---------------------------------------------------------------------------
# 3) Insert two additional blank paragraphs (same "Retraitcorpsdetexte"
#    style as the blank paragraph that already precedes it) right before the
#    "Pratique" Titre1 heading.
# ---------------------------------------------------------------------------
$pratiqueTitle = Find-ParaByStyleLike "Heading 1" "Pratique*"
if ($pratiqueTitle -eq $null) {
    Write-Output "ERROR: could not locate Pratique heading"
} else {
    $blankBefore = $pratiqueTitle.Previous(1)
    $blankBefore.Range.InsertParagraphAfter()
    $blankBefore.Range.InsertParagraphAfter()
}

Write-Output "edit complete"
